$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches the source inlineStr cells)
# without leaving a quote-prefix / text-numFmt style behind on the cell.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "308.59"
Set-TextValue $ws.Range("E2") "6.94%"
Set-TextValue $ws.Range("D3") "32.22"
Set-TextValue $ws.Range("E3") "10.33%"
Set-TextValue $ws.Range("D4") "5.324"
Set-TextValue $ws.Range("E4") "4.88%"
Set-TextValue $ws.Range("D5") "0.07472"
Set-TextValue $ws.Range("E5") "11.89%"
Set-TextValue $ws.Range("D6") "7.798"
Set-TextValue $ws.Range("E6") "6.02%"
Set-TextValue $ws.Range("D7") "3.667"
Set-TextValue $ws.Range("E7") "7.62%"
Set-TextValue $ws.Range("D8") "1.572"
Set-TextValue $ws.Range("E8") "16.61%"
Set-TextValue $ws.Range("D9") "0.9103"
Set-TextValue $ws.Range("E9") "-0.86%"
Set-TextValue $ws.Range("D10") "0.01712"
Set-TextValue $ws.Range("E10") "2,553.70%"
Set-TextValue $ws.Range("D11") "0.1678"
Set-TextValue $ws.Range("E11") "6.14%"
Set-TextValue $ws.Range("D12") "0.07577"
Set-TextValue $ws.Range("E12") "11.55%"
Set-TextValue $ws.Range("D13") "0.08044"
Set-TextValue $ws.Range("E13") "4.88%"
Set-TextValue $ws.Range("D14") "0.03033"
Set-TextValue $ws.Range("E14") "3.47%"
Set-TextValue $ws.Range("D15") "0.09862"
Set-TextValue $ws.Range("E15") "9.70%"
Set-TextValue $ws.Range("D16") "0.001519"
Set-TextValue $ws.Range("E16") "-4.09%"
Set-TextValue $ws.Range("D17") "0.04546"
Set-TextValue $ws.Range("E17") "0.81%"
Set-TextValue $ws.Range("D18") "0.006485"
Set-TextValue $ws.Range("E18") "3.44%"
Set-TextValue $ws.Range("D19") "3.489"
Set-TextValue $ws.Range("E19") "1.02%"
Set-TextValue $ws.Range("D20") "2.241"
Set-TextValue $ws.Range("E20") "0.97%"
Set-TextValue $ws.Range("D21") "0.3266"
Set-TextValue $ws.Range("E21") "1.66%"
Set-TextValue $ws.Range("D23") "4.184"
Set-TextValue $ws.Range("E23") "2.77%"
Set-TextValue $ws.Range("D24") "0.1618"
Set-TextValue $ws.Range("E24") "2.25%"
Set-TextValue $ws.Range("D25") "0.001213"
Set-TextValue $ws.Range("E25") "1.93%"
Set-TextValue $ws.Range("D26") "0.004504"
Set-TextValue $ws.Range("E26") "9.25%"
Set-TextValue $ws.Range("D27") "0.0001299"
Set-TextValue $ws.Range("E27") "8.32%"
Set-TextValue $ws.Range("D28") "0.0001739"
Set-TextValue $ws.Range("E28") "7.51%"
Set-TextValue $ws.Range("D40") "0.04530"
Set-TextValue $ws.Range("E40") "7.90%"
Set-TextValue $ws.Range("D41") "0.007179"
Set-TextValue $ws.Range("E41") "7.04%"
Set-TextValue $ws.Range("D42") "0.1363"
Set-TextValue $ws.Range("E42") "9.99%"
Set-TextValue $ws.Range("D43") "0.002258"
Set-TextValue $ws.Range("E43") "8.06%"
Set-TextValue $ws.Range("D44") "0.01392"
Set-TextValue $ws.Range("E44") "4.55%"
Set-TextValue $ws.Range("D45") "0.00006147"
Set-TextValue $ws.Range("E45") "7.34%"
Set-TextValue $ws.Range("D47") "0.01298"
Set-TextValue $ws.Range("E47") "-0.67%"

# Row 22: only Volume(1h) changed, Price stayed the same
Set-TextValue $ws.Range("E22") "1.95%"

# Row 46: only Price changed, Volume(1h) stayed the same
Set-TextValue $ws.Range("D46") "1.892"

